{"js": "// Center images in book-type files: paragraphs that consist solely of an\n// inline picture (the book cover / avatar / logo images) get the \"Compact\"\n// paragraph style applied, matching how the rest of the document already\n// marks tightly-spaced content.\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// Collect the inline-picture collection for every paragraph first so we can\n// batch the loads in one round-trip.\nconst pictureCollections = paragraphs.items.map((p) => p.inlinePictures);\npictureCollections.forEach((pics) => pics.load(\"items\"));\nawait context.sync();\n\n// Apply the \"Compact\" paragraph style to every paragraph that contains at\n// least one inline picture.\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (pictureCollections[i].items.length > 0) {\n    paragraphs.items[i].style = \"Compact\";\n  }\n}\nawait context.sync();\n", "ps1": "# Center images in book-type files: any paragraph that consists of an\n# inline picture (book cover / avatar / logo images) gets the \"Compact\"\n# paragraph style applied, matching how the rest of the document already\n# marks tightly-spaced content.\n$d = $word.ActiveDocument\n\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $para = $d.Paragraphs.Item($i)\n    if ($para.Range.InlineShapes.Count -gt 0) {\n        $para.Style = \"Compact\"\n    }\n}\n"}
